# Edit Book2.xlsx / Sheet1 per the target diff:
#  - Update a few "Requests" text values in column H
#  - Set column B (the "X" checkbox column) to 1 for rows 3-9
#  - Move the active selection to B9
#  - Row heights for rows 4-9 grow slightly (a side effect of Excel
#    auto-adjusting row height once the bold/arial "X" entries are added)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update text values -----------------------------------------------
$ws.Range("H5").Value = "Only After 5:00pm"
$ws.Range("H6").Value = "no 1/12/14-1/20/14, no 1/20/14"
$ws.Range("H7").Value = "No 1/17/14"

# --- Mark teams with an X in column B for rows 3-9 ---------------------
for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}

# --- Row heights for rows 4-9 grow slightly -----------------------------
$ws.Rows.Item(4).RowHeight = 14.05
$ws.Rows.Item(5).RowHeight = 14.05
$ws.Rows.Item(6).RowHeight = 14.05
$ws.Rows.Item(7).RowHeight = 14.05
$ws.Rows.Item(8).RowHeight = 14.05
$ws.Rows.Item(9).RowHeight = 14.9

# --- Move the active cell/selection to B9 -------------------------------
$ws.Activate()
$ws.Range("B9").Select()
